$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (leg max ROM columns) updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON, meanEMG) updated values
$ws.Range("B2").Value = 94.77317618170018
$ws.Range("C2").Value = 93.391657235873353
$ws.Range("D2").Value = 93.19621952464972
$ws.Range("E2").Value = 94.440979878641414

# Row 3 (STR, meanEMG) updated values
$ws.Range("B3").Value = 93.934020850837427
# C3 is cleared (no longer has data in the updated dataset)
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 91.621125973895516
$ws.Range("E3").Value = 95.128685951079902

# Reflect the updated active selection over the edited data range
[void]$ws.Range("B1:E3").Select()
